# Applies the "Modified costs data" change:
#  - Adds VarOM values for Gas-CCS-95 / Gas-CCS-97 (G9, G10)
#  - Adds three new BESS technology rows (11, 12, 13) with the usual
#    CAPEX/PMT and CAPEX*1000 formulas copied down from row 10
#  - Changes the number format of style index 3 (used by the "Ramping"
#    column, K) from General to 0.00%, and reassigns that style to the
#    new BESS Capacity Factor cells in column E instead of K2
#  - Updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the two missing VarOM values in the existing CCS rows ---
$ws.Range("G9").Value = 1.6
$ws.Range("G10").Value = 1.6

# --- New data rows for the BESS technologies ---
$ws.Range("A11").Value = "BESS-2H"
$ws.Range("B11").Value = 943
$ws.Range("C11").Value = 24
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.083
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("K11").Value = 1

$ws.Range("A12").Value = "BESS-6H"
$ws.Range("B12").Value = 2321
$ws.Range("C12").Value = 56
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.25
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("K12").Value = 1

$ws.Range("A13").Value = "BESS-10H"
$ws.Range("B13").Value = 3250
$ws.Range("C13").Value = 88
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.417
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("K13").Value = 1

# Copy the number formatting from row 10 (CAPEX $MW-y / OM $/MW-y columns)
# down onto the new rows before writing their formulas, so the new cells
# reuse the same cellXf (currency style) instead of Excel fabricating a
# brand-new style.
$ws.Range("I10").Copy()
$ws.Range("I11:I13").PasteSpecial(-4122)

$ws.Range("I11").Formula = "=PMT(5%,30,-B11*1000)"
$ws.Range("I12").Formula = "=PMT(5%,30,-B12*1000)"
$ws.Range("I13").Formula = "=PMT(5%,30,-B13*1000)"

$ws.Range("J11").Formula = "=C11*1000"
$ws.Range("J12").Formula = "=C12*1000"
$ws.Range("J13").Formula = "=C13*1000"

# --- Capacity Factor formatting for the new rows ---
# E12 uses the plain 0% style already used elsewhere in the column.
$ws.Range("E12").NumberFormat = "0%"

# E11 / E13 reuse the style currently sitting at K2 (cellXf index 3),
# but change its number format from General to 0.00%.
$ws.Range("K2").NumberFormat = "0.00%"
$ws.Range("E11").NumberFormat = "0.00%"
$ws.Range("E13").NumberFormat = "0.00%"

# K2 itself goes back to the default (General) style now that cellXf 3
# has been repurposed for the percentage formatting above.
$ws.Range("K2").ClearFormats()

$excel.Calculate()

# --- Update the active cell / selection like the author's session ---
$ws.Range("G18").Select()
